$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AnnotatorTestSelenium")

# Update row 2 data: use a different cadd resource
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 158796
$ws.Range("D2").Value = "A"
$ws.Range("E2").Value = "C"

# Move the active selection from C4 to C3
$ws.Activate()
$ws.Range("C3").Select()
